# Apply edits described by the commit diff:
# 1. Rename sheet "Sheet1" -> "corrupt"
# 2. Update the absPath (best-effort; may not be exposed via the object model)
# 3. Move sheet view: topLeftCell B1 -> C1, selection E14 -> C13
# 4. Add a thin left border to a new style, applied to H8:H11 with a single-space value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename sheet
$ws.Name = "corrupt"

# 3. Update view: scroll position + selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3   # C is column 3 -> topLeftCell C1
$ws.Range("C13").Select()

# 4. Add single-space values with a left-thin-bordered style to H8:H11
$cells = @("H8", "H9", "H10", "H11")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value = " "
    $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft = 7, xlContinuous = 1
    $rng.Borders.Item(7).Weight = 2      # xlThin = 2
}
